# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet per the "Updated symbol list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    # Force text storage so numeric/percent-looking strings are not
    # reinterpreted as numbers (matches the source file's inlineStr cells).
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "301.50"
Set-TextValue $ws "E2" "-0.69%"
Set-TextValue $ws "D3" "31.41"
Set-TextValue $ws "E3" "-1.79%"
Set-TextValue $ws "D4" "5.119"
Set-TextValue $ws "E4" "-2.43%"
Set-TextValue $ws "D5" "0.07349"
Set-TextValue $ws "E5" "-2.60%"
Set-TextValue $ws "D6" "2.249"
Set-TextValue $ws "E6" "48.53%"
Set-TextValue $ws "D7" "7.933"
Set-TextValue $ws "E7" "0.05%"
Set-TextValue $ws "D8" "3.789"
Set-TextValue $ws "E8" "-0.70%"
Set-TextValue $ws "D9" "0.9182"
Set-TextValue $ws "E9" "-0.63%"
Set-TextValue $ws "D10" "0.1701"
Set-TextValue $ws "E10" "0.12%"
Set-TextValue $ws "D11" "0.07547"
Set-TextValue $ws "E11" "-4.62%"
Set-TextValue $ws "D12" "0.08183"
Set-TextValue $ws "E12" "1.68%"
Set-TextValue $ws "D13" "0.03027"
Set-TextValue $ws "E13" "-0.77%"
Set-TextValue $ws "D14" "0.09930"
Set-TextValue $ws "E14" "0.14%"
Set-TextValue $ws "D15" "0.001494"
Set-TextValue $ws "E15" "-1.26%"
Set-TextValue $ws "D16" "0.006133"
Set-TextValue $ws "E16" "-3.27%"
Set-TextValue $ws "D17" "3.462"
Set-TextValue $ws "E17" "0.41%"
Set-TextValue $ws "D18" "2.219"
Set-TextValue $ws "E18" "-0.53%"
Set-TextValue $ws "D19" "0.3305"
Set-TextValue $ws "E19" "0.18%"
Set-TextValue $ws "D20" "0.1340"
Set-TextValue $ws "E20" "0.99%"
Set-TextValue $ws "D21" "4.656"
Set-TextValue $ws "D22" "0.04659"
Set-TextValue $ws "E22" "1.37%"
Set-TextValue $ws "D23" "0.1569"
Set-TextValue $ws "E23" "-3.02%"
Set-TextValue $ws "D24" "0.001228"
Set-TextValue $ws "E24" "0.93%"
Set-TextValue $ws "D25" "0.004467"
Set-TextValue $ws "E25" "0.13%"
Set-TextValue $ws "D26" "0.0001301"
Set-TextValue $ws "E26" "-7.00%"
Set-TextValue $ws "E27" "49.51%"
Set-TextValue $ws "D39" "0.01728"
Set-TextValue $ws "E39" "1.93%"
Set-TextValue $ws "D40" "0.04523"
Set-TextValue $ws "E40" "0.67%"
Set-TextValue $ws "D41" "0.007237"
Set-TextValue $ws "D42" "0.1344"
Set-TextValue $ws "E42" "-0.67%"
Set-TextValue $ws "D43" "0.002231"
Set-TextValue $ws "E43" "7.37%"
Set-TextValue $ws "D44" "0.01075"
Set-TextValue $ws "E44" "-22.07%"
Set-TextValue $ws "D45" "0.00006307"
Set-TextValue $ws "E45" "2.12%"
Set-TextValue $ws "E46" "-23.07%"
Set-TextValue $ws "E47" "-55.49%"
